# Flip RunMode flags from "Y" to "N" on the data-driven test sheets
# (Testcases, the CheckItems data grid, and LoginTest), and leave the
# workbook positioned/selected the way the author last left it: LoginTest
# tab active with F5 selected.

$wb = $excel.ActiveWorkbook

# --- Testcases sheet: RunMode column (C) ---
$wsTestcases = $wb.Worksheets.Item("Testcases")
$wsTestcases.Range("C3").Value = "N"
$wsTestcases.Range("C4").Value = "N"
$wsTestcases.Range("C5").Value = "N"

# --- CheckItems sheet: RunMode column (C) in the data grid ---
$wsCheckItems = $wb.Worksheets.Item("CheckItems")
$wsCheckItems.Range("C3").Value = "N"
$wsCheckItems.Range("C5").Value = "N"

# --- LoginTest sheet: RunMode column (F) ---
$wsLoginTest = $wb.Worksheets.Item("LoginTest")
$wsLoginTest.Range("F3").Value = "N"

# --- Restore per-sheet selections as last left by the author ---
$wsTestcases.Activate() | Out-Null
$wsTestcases.Range("C2").Select() | Out-Null

$wsTeststeps = $wb.Worksheets.Item("Teststeps")
$wsTeststeps.Activate() | Out-Null
$wsTeststeps.Range("C14").Select() | Out-Null

$wsCheckItems.Activate() | Out-Null
$wsCheckItems.Range("C6").Select() | Out-Null

$wsLoginTest.Activate() | Out-Null
$wsLoginTest.Range("F5").Select() | Out-Null
